$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# maximum_investment_capacity_per_year: 10000 -> 1000000, formatted in scientific notation
$ws.Range("B13").Value = 1000000
$ws.Range("B13").NumberFormat = "0.00E+00"

# realistic_candidate_capacities: TRUE -> FALSE
$ws.Range("B16").Value = $false

# realistic_candidate_capacities_for_future: TRUE -> FALSE
$ws.Range("B17").Value = $false

# Move selection to B19 (just below the data) as saved by the author
[void]$ws.Range("B19").Select()
